{"js": "// Merge the split runs of the Title paragraph and the Abstract paragraph\n// into single runs, without changing the visible text, by replacing each\n// paragraph's text in place (Office.js collapses a paragraph's text into\n// a single run when it is rewritten this way).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  const style = paragraph.style;\n  const text = paragraph.text;\n\n  if (style === \"Title\" && text === \"Answers: Trigonometry (degrees)\") {\n    paragraph.insertText(\"Answers: Trigonometry (degrees)\", \"Replace\");\n  } else if (\n    style === \"Abstract\" &&\n    text === \"Answers to the questions on trigonometry, using degrees to measure angles.\"\n  ) {\n    paragraph.insertText(\n      \"Answers to the questions on trigonometry, using degrees to measure angles.\",\n      \"Replace\"\n    );\n  }\n}\n\nawait context.sync();\n", "ps1": "# Merge the split runs that make up the Title paragraph and the Abstract\n# paragraph into a single run each, without changing the visible text.\n# Using Find/Replace on the exact paragraph text re-writes the matched\n# span as one run, which is how Word collapses multiple runs that carry\n# identical (default) formatting into one when the text is replaced.\n\n$d = $word.ActiveDocument\n\nfunction Merge-ParagraphRuns($searchText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $searchText\n    $find.Replacement.Text = $searchText\n    [void]$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n\nMerge-ParagraphRuns \"Answers: Trigonometry (degrees)\"\nMerge-ParagraphRuns \"Answers to the questions on trigonometry, using degrees to measure angles.\"\n"}
